# Edit script: applies the betexplorer CFL group B update described by the diff.
# - Swaps / rotates the F:V (match data) contents of several row groups
#   (the A:E "key" columns - Indice, pais, torneio, temporada, data_partida -
#   stay untouched, only the match-result/odds columns move).
# - Appends four brand-new match rows (105-108) at the bottom, copying the
#   formatting of the last existing row and then filling in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $range1 = $ws.Range("F$($r1):V$($r1)")
    $range2 = $ws.Range("F$($r2):V$($r2)")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value = $v2
    $range2.Value = $v1
}

# --- Group 1: rows 5 <-> 6 -------------------------------------------------
Swap-Rows 5 6

# --- Group 2: rows 70, 71, 72 rotate backward ------------------------------
# new70 = old72 ; new71 = old70 ; new72 = old71
$v70 = $ws.Range("F70:V70").Value2
$v71 = $ws.Range("F71:V71").Value2
$v72 = $ws.Range("F72:V72").Value2
$ws.Range("F70:V70").Value = $v72
$ws.Range("F71:V71").Value = $v70
$ws.Range("F72:V72").Value = $v71

# --- Group 3: rows 79, 80, 81 rotate backward ------------------------------
# new79 = old81 ; new80 = old79 ; new81 = old80
$v79 = $ws.Range("F79:V79").Value2
$v80 = $ws.Range("F80:V80").Value2
$v81 = $ws.Range("F81:V81").Value2
$ws.Range("F79:V79").Value = $v81
$ws.Range("F80:V80").Value = $v79
$ws.Range("F81:V81").Value = $v80

# --- Group 4: rows 83 <-> 84 ------------------------------------------------
Swap-Rows 83 84

# --- Group 5: rows 94, 95, 96 rotate forward --------------------------------
# new94 = old95 ; new95 = old96 ; new96 = old94
$v94 = $ws.Range("F94:V94").Value2
$v95 = $ws.Range("F95:V95").Value2
$v96 = $ws.Range("F96:V96").Value2
$ws.Range("F94:V94").Value = $v95
$ws.Range("F95:V95").Value = $v96
$ws.Range("F96:V96").Value = $v94

# --- Group 6: rows 103 <-> 104 ----------------------------------------------
Swap-Rows 103 104

# --- Group 7: append new rows 105-108 ---------------------------------------
# Copy formatting from the current last row (104) down into the new rows,
# then populate the values.
$ws.Range("A104:V104").Copy() | Out-Null
$ws.Range("A105:V108").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row=105; A=104; E=45234.42708333334; F="Jablonec B"; G=1; H="Prepere"; I=1;
       J=1.41; K="04/11/2023 00:42"; L=1.6;  M="04/11/2023 10:13";
       N=4.7;  O="04/11/2023 00:42"; P=4.22; Q="04/11/2023 10:13";
       R=5.58; S="04/11/2023 00:42"; T=4.45; U="04/11/2023 10:12";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-b/jablonec-prepere/82FdtFff/" },
    @{ Row=106; A=105; E=45234.4375; F="Kolin"; G=1; H="Pardubice B"; I=1;
       J=1.31; K="04/11/2023 01:13"; L=1.39; M="04/11/2023 10:00";
       N=5.28; O="04/11/2023 01:13"; P=5.39; Q="04/11/2023 10:11";
       R=6.62; S="04/11/2023 01:13"; T=5.59; U="04/11/2023 10:11";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-b/kolin-pardubice/bi89wDvD/" },
    @{ Row=107; A=106; E=45234.58333333334; F="Banik Most-Sous"; G=2; H="Zapy"; I=2;
       J=2.33; K="04/11/2023 12:13"; L=2.16; M="04/11/2023 13:48";
       N=3.56; O="04/11/2023 12:13"; P=3.16; Q="04/11/2023 13:48";
       R=2.61; S="04/11/2023 12:13"; T=3.27; U="04/11/2023 13:48";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-b/banik-most-sous-zapy/Qc0e1GWJ/" },
    @{ Row=108; A=107; E=45234.58333333334; F="Zivanice"; G=3; H="Usti nad Labem"; I=1;
       J=3.72; K="04/11/2023 12:13"; L=3.59; M="04/11/2023 13:59";
       N=3.75; O="04/11/2023 12:13"; P=3.98; Q="04/11/2023 13:59";
       R=1.79; S="04/11/2023 12:13"; T=1.8;  U="04/11/2023 13:59";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-b/zivanice-usti-nad-labem/xr7DxXgJ/" }
)

foreach ($rd in $newRows) {
    $r = $rd.Row
    $ws.Range("A$r").Value = $rd.A
    $ws.Range("B$r").Value = "czech-republic"
    $ws.Range("C$r").Value = "cfl-group-b"
    $ws.Range("D$r").Value = "2023-2024"
    $ws.Range("E$r").Value = $rd.E
    $ws.Range("F$r").Value = $rd.F
    $ws.Range("G$r").Value = $rd.G
    $ws.Range("H$r").Value = $rd.H
    $ws.Range("I$r").Value = $rd.I
    $ws.Range("J$r").Value = $rd.J
    $ws.Range("K$r").Value = $rd.K
    $ws.Range("L$r").Value = $rd.L
    $ws.Range("M$r").Value = $rd.M
    $ws.Range("N$r").Value = $rd.N
    $ws.Range("O$r").Value = $rd.O
    $ws.Range("P$r").Value = $rd.P
    $ws.Range("Q$r").Value = $rd.Q
    $ws.Range("R$r").Value = $rd.R
    $ws.Range("S$r").Value = $rd.S
    $ws.Range("T$r").Value = $rd.T
    $ws.Range("U$r").Value = $rd.U
    $ws.Range("V$r").Value = $rd.V
}
